$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -8
$ws.Range("F3").Value = -2
$ws.Range("F7").Value = -2
$ws.Range("F8").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("F18").Value = -2
